# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 21 de Agosto de 2020 a las 09:45"

# Row 6 - India
$ws.Range("B6").Value = 2910032
$ws.Range("C6").Value = 5703
$ws.Range("D6").Value = 2160059
$ws.Range("E6").Value = 694971
$ws.Range("G6").Value = 27
$ws.Range("H6").Value = 55002

# Row 7 - Rusia
$ws.Range("B7").Value = 946976
$ws.Range("C7").Value = 4870
$ws.Range("D7").Value = 761330
$ws.Range("E7").Value = 169457
$ws.Range("G7").Value = 90
$ws.Range("H7").Value = 16189

# Row 57 - Armenia
$ws.Range("B57").Value = 42477
$ws.Range("C57").Value = 158
$ws.Range("D57").Value = 35693
$ws.Range("E57").Value = 5942
$ws.Range("G57").Value = 6
$ws.Range("H57").Value = 842

# Row 60 - Afganistan
$ws.Range("B60").Value = 37894
$ws.Range("C60").Value = 38
$ws.Range("D60").Value = 28016
$ws.Range("E60").Value = 8493

# Row 108 - Hungria
$ws.Range("B108").Value = 5098
$ws.Range("C108").Value = 52
$ws.Range("D108").Value = 3681
$ws.Range("E108").Value = 806
$ws.Range("G108").Value = 2
$ws.Range("H108").Value = 611

# Row 133 - Estonia
$ws.Range("B133").Value = 2244
$ws.Range("C133").Value = 17
$ws.Range("D133").Value = 2011
$ws.Range("E133").Value = 170

# Row 149 - Georgia
$ws.Range("B149").Value = 1385
$ws.Range("C149").Value = 15
$ws.Range("D149").Value = 1128
$ws.Range("E149").Value = 240

# Row 150 - Letonia
$ws.Range("B150").Value = 1330
$ws.Range("C150").Value = 3
$ws.Range("E150").Value = 204

# Row 206 - Laos
$ws.Range("D206").Value = 20
$ws.Range("E206").Value = 2

$wb.Save()
